# Daily attendance processing - 2026-01-21 23:38:20
#
# In the "Recorded By" column (G), normalize the ordering of the two
# recorder names for entries that were recorded by a named user together
# with "System": the literal tag "System" should be listed first,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# and  "admin@admin.com, System"    -> "System, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($value -eq "admin@admin.com, System") {
        $cell.Value2 = "System, admin@admin.com"
    }
}
